# AI-generated content refresh + trimming of the last two slides.
$p = $ppt.ActivePresentation

# Helper: replace a paragraph's text with a clean, single run. Clearing the
# paragraph first avoids the host's "keep common prefix/suffix run" merge
# logic, which otherwise splits the replacement into several <a:r> runs.
function Set-ParaText($paragraph, [string]$text) {
    $paragraph.Text = ""
    $paragraph.Text = $text
}

# --- Remove the final two slides (old "The Future of Quantum Computing" and
#     "Conclusion..." slides) -- their content / role has been folded into
#     the updated slide 10. Delete from the back so indices stay valid. ---
$p.Slides.Item(12).Delete()
$p.Slides.Item(11).Delete()

# --- Slide 2: Title slide ---
$s2 = $p.Slides.Item(2)
Set-ParaText $s2.Shapes.Item(1).TextFrame.TextRange.Paragraphs(1,1) "Quantum Computing: A New Era of Computation"
Set-ParaText $s2.Shapes.Item(2).TextFrame.TextRange.Paragraphs(2,1) "Created by: guppi"

# --- Slide 3: Index ---
$s3 = $p.Slides.Item(3)
Set-ParaText $s3.Shapes.Item(2).TextFrame.TextRange.Paragraphs(2,1) "1.Introduction to Quantum Computing2.Quantum Bits (Qubits) and Superposition3.Quantum Entanglement and its Implications4.Quantum Algorithms and their Advantages5.Current Applications and Future Potential6.Challenges and Limitations of Quantum Computing7.The Future of Quantum ComputingConclusion"

# --- Slide 4: Introduction to Quantum Computing ---
$s4 = $p.Slides.Item(4)
Set-ParaText $s4.Shapes.Item(2).TextFrame.TextRange.Paragraphs(2,1) "Classical computers use bits representing 0 or 1.Quantum computers use qubits, leveraging superposition and entanglement.Quantum computing harnesses quantum mechanics to solve complex problems.Potential to revolutionize various fields like medicine, materials science, and finance.Significant advancements in hardware and software are ongoing.Explores computation beyond the limitations of classical computers.Focuses on solving problems intractable for classical systems."

# --- Slide 5: Quantum Bits (Qubits) ---
$s5 = $p.Slides.Item(5)
Set-ParaText $s5.Shapes.Item(1).TextFrame.TextRange.Paragraphs(1,1) "Quantum Bits (Qubits) and Superposition"
Set-ParaText $s5.Shapes.Item(2).TextFrame.TextRange.Paragraphs(2,1) "Qubits can represent 0, 1, or a combination of both simultaneously (superposition).Superposition allows quantum computers to explore multiple possibilities concurrently.Different physical systems can be used to represent qubits (e.g., trapped ions, superconducting circuits).Measurement collapses the superposition into a definite 0 or 1.Control and manipulation of qubits are crucial for quantum computation.Coherence time (how long a qubit maintains superposition) is a key challenge."

# --- Slide 6: Quantum Entanglement ---
$s6 = $p.Slides.Item(6)
Set-ParaText $s6.Shapes.Item(1).TextFrame.TextRange.Paragraphs(1,1) "Quantum Entanglement and its Implications"
Set-ParaText $s6.Shapes.Item(2).TextFrame.TextRange.Paragraphs(2,1) "Entanglement links two or more qubits, regardless of distance.Measuring the state of one entangled qubit instantly reveals the state of the others.Entanglement enables powerful quantum algorithms and computations.Einstein called it Understanding and harnessing entanglement is vital for quantum technologies.Entanglement is a key resource for quantum communication and cryptography."

# --- Slide 7: Quantum Algorithms ---
$s7 = $p.Slides.Item(7)
Set-ParaText $s7.Shapes.Item(1).TextFrame.TextRange.Paragraphs(1,1) "Quantum Algorithms and their Advantages"
Set-ParaText $s7.Shapes.Item(2).TextFrame.TextRange.Paragraphs(2,1) "Shor's algorithm: efficiently factors large numbers (cryptography implications).Grover's algorithm: speeds up database searches quadratically.Quantum algorithms offer exponential speedups over classical algorithms for specific problems.Development of new quantum algorithms is an active area of research.Quantum machine learning algorithms are emerging.Quantum simulation promises breakthroughs in materials science and drug discovery."

# --- Slide 8: Applications ---
$s8 = $p.Slides.Item(8)
Set-ParaText $s8.Shapes.Item(1).TextFrame.TextRange.Paragraphs(1,1) "Current Applications and Future Potential"
Set-ParaText $s8.Shapes.Item(2).TextFrame.TextRange.Paragraphs(2,1) "Drug discovery and materials science: simulating molecular interactions.Financial modeling: optimizing portfolios and risk management.Cryptography: developing quantum-resistant encryption methods.Optimization problems: solving complex logistics and scheduling tasks.Artificial intelligence: enhancing machine learning algorithms.Quantum sensing and metrology: improving precision measurements."

# --- Slide 9: Challenges ---
$s9 = $p.Slides.Item(9)
Set-ParaText $s9.Shapes.Item(1).TextFrame.TextRange.Paragraphs(1,1) "Challenges and Limitations of Quantum Computing"
Set-ParaText $s9.Shapes.Item(2).TextFrame.TextRange.Paragraphs(2,1) "Building and maintaining stable qubits is technologically challenging.Error correction is crucial due to qubit decoherence.Quantum computers are currently expensive and require specialized environments.Scaling up the number of qubits while maintaining coherence is a major hurdle.Developing quantum algorithms requires specialized expertise.Limited availability and accessibility of quantum computing resources."

# --- Slide 10: was "Prominent Players...", now "The Future of Quantum Computing" ---
$s10 = $p.Slides.Item(10)
Set-ParaText $s10.Shapes.Item(1).TextFrame.TextRange.Paragraphs(1,1) "The Future of Quantum Computing"
Set-ParaText $s10.Shapes.Item(2).TextFrame.TextRange.Paragraphs(2,1) "Continued advancements in hardware and software are expected.Increased accessibility and affordability of quantum computing resources.Collaboration between academia, industry, and government is essential.Potential for disruptive innovations across various sectors.Ethical considerations and societal impacts need careful consideration.Quantum computing will likely coexist with classical computing, complementing its strengths.A new era of scientific discovery and technological advancement is on the horizon."
